$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve default (unstyled) cell style to reapply after forcing text format,
# so numeric-looking price strings stay text without leaving stray styles behind.
$defaultStyle = $ws.Range("B2").Style

# Force column D to Text format so values like "0.393" are stored as strings,
# matching the workbook convention of inline/shared-string prices.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "92.410.47"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "3.106.11"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "241.14"
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").Value = "615.43"
$ws.Range("E6").Value = "  -1.62%  "

$ws.Range("E7").Value = "  -6.48%  "

$ws.Range("D8").Value = "0.393"
$ws.Range("E8").Value = "  +4.20%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").Value = "3.099.85"
$ws.Range("E10").Value = "  -1.53%  "

$ws.Range("D11").Value = "0.730"
$ws.Range("E11").Value = "  -4.42%  "

$ws.Range("E12").Value = "  -1.70%  "

$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("D14").Value = "34.39"
$ws.Range("E14").Value = "  -3.59%  "

$ws.Range("D15").Value = "91.921.22"
$ws.Range("E15").Value = "  -0.13%  "

$ws.Range("E16").Value = "  -0.63%  "

$ws.Range("D17").Value = "3.677.05"
$ws.Range("E17").Value = "  -1.06%  "

$ws.Range("D18").Value = "3.089.28"
$ws.Range("E18").Value = "  -0.64%  "

$ws.Range("E19").Value = "  -2.85%  "

$ws.Range("D20").Value = "14.74"
$ws.Range("E20").Value = "  -2.61%  "

$ws.Range("E21").Value = "  -0.86%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "9.37"
$ws.Range("E22").Value = "  +1.40%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "447.64"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("E24").Value = "  -6.13%  "

$ws.Range("D25").Value = "5.61"
$ws.Range("E25").Value = "  -2.15%  "

$ws.Range("D26").Value = "87.12"
$ws.Range("E26").Value = "  -4.03%  "

$ws.Range("D27").Value = "11.74"
$ws.Range("E27").Value = "  -2.55%  "

$ws.Range("D28").Value = "3.265.72"

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("E30").Value = "  +9.65%  "

$ws.Range("E31").Value = "  -5.24%  "

$ws.Range("D32").Value = "0.168"
$ws.Range("E32").Value = "  -4.67%  "

$ws.Range("D33").Value = "9.22"
$ws.Range("E33").Value = "  -1.53%  "

$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +5.75%  "

$ws.Range("D35").Value = "8.04"
$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").Value = "0.161"
$ws.Range("E36").Value = "  -5.98%  "

$ws.Range("D37").Value = "4.25"
$ws.Range("E37").Value = "  -2.69%  "

$ws.Range("D38").Value = "26.20"
$ws.Range("E38").Value = "  -2.06%  "

$ws.Range("E39").Value = "  -1.00%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "1.30"
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "480.29"
$ws.Range("E41").Value = "  -4.51%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "3.45"
$ws.Range("E42").Value = "  -3.81%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "0.434"
$ws.Range("E43").Value = "  +1.40%  "

$ws.Range("D44").Value = "22.87"
$ws.Range("E44").Value = "  +2.71%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").Value = "159.37"
$ws.Range("E46").Value = "  +3.20%  "

$ws.Range("D47").Value = "1.90"
$ws.Range("E47").Value = "  -2.97%  "

$ws.Range("D48").Value = "0.694"
$ws.Range("E48").Value = "  -1.93%  "

$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("D50").Value = "0.0333"
$ws.Range("E50").Value = "  +3.45%  "

$ws.Range("D51").Value = "44.14"
$ws.Range("E51").Value = "  -0.82%  "

# Restore the original default style to column D (keeps text type, drops the
# temporary Text number-format styling so cells match the source workbook).
$ws.Range("D2:D51").Style = $defaultStyle
